$p = $ppt.ActivePresentation

# ----------------------------------------------------------------
# Slide 12 : "OUR SCHEDULE" slide
# ----------------------------------------------------------------
$s12 = $p.Slides.Item(12)

# Shape 3 ("TextBox 4") holds the bulleted schedule list.
$scheduleBox = $s12.Shapes.Item(3)
# Suspend autofit while editing the runs so the shape keeps its
# original size (matches the target OOXML, which leaves <a:ext> alone).
$scheduleBox.TextFrame.AutoSize = 0
$scheduleTr = $scheduleBox.TextFrame.TextRange
$scheduleTr.Paragraphs(2).Runs(1).Text = " Introduction & Environment : 4 hours – 1 sessions."
$scheduleTr.Paragraphs(3).Runs(1).Text = " HTML & CSS & Bootstrap 4 : 32 hours – 8 sessions."
$scheduleTr.Paragraphs(4).Runs(1).Text = "JavaScript & jQuery & AJAX & JSON : 16 hours – 4 sessions."
$scheduleTr.Paragraphs(5).Runs(1).Text = "Project News Website Template : 8 hours – 2 sessions."
$scheduleBox.TextFrame.AutoSize = 1
$scheduleBox.Height = 288.70216098425203

# Shape 4 ("TextBox 6") holds the footnote; it also gets wider.
$footnoteBox = $s12.Shapes.Item(4)
$footnoteBox.TextFrame.AutoSize = 0
$footnoteTr = $footnoteBox.TextFrame.TextRange
$footnoteTr.Paragraphs(1).Runs(1).Text = "* Git Hub, Git Page, VS Code, Command Line will be guided in some session."
$footnoteBox.TextFrame.AutoSize = 1
$footnoteBox.Width = 595.8375253149607
$footnoteBox.Height = 35.68011374015748

# ----------------------------------------------------------------
# Slide 3 : "HELLO, I'M NAM." bio slide
# ----------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$bioBox = $s3.Shapes.Item(3)
$bioBox.TextFrame.AutoSize = 0
$bioTr = $bioBox.TextFrame.TextRange
$bioTr.Paragraphs(1).Runs(1).Text = "Software Engineering at Axon Active Viet Nam."
$bioBox.TextFrame.AutoSize = 1
$bioBox.Height = 230.5396412992126

